$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A
$ws.Columns.Item(1).ColumnWidth = 16.1

# Update "seq lcc2" -> "seq local -O3" in A6
$ws.Range("A6").Value = "seq local -O3"

# New row A7 label
$ws.Range("A7").Value = "omp "

# Fill in the new benchmark numbers for row 6 (seq local -O3)
$ws.Range("B6").Value = 0.000172
$ws.Range("C6").Value = 0.002159
$ws.Range("D6").Value = 0.260026
$ws.Range("E6").Value = 4.422902
$ws.Range("F6").Value = 88.539109
$ws.Range("G6").Value = 147.002965

# Bump the precision of the Menlo-styled measurement row to 6 decimals and
# propagate that format to the newly-filled row 6
$ws.Range("B5:G5").NumberFormat = "0.000000"
$ws.Range("B5:G5").Copy()
$ws.Range("B6:G6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Select D15 (matches the saved view state)
$ws.Range("D15").Select()

# Set up the printable page (A4 portrait)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
